$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6
$ws.Range("AE6").Value = 26
$ws.Range("AF6").Value = 101
$ws.Range("AI6").Value = 29
$ws.Range("AJ6").Value = 21
$ws.Range("AK6").Value = 67
$ws.Range("AL6").Value = 51
$ws.Range("AM6").Value = 67
$ws.Range("AN6").Value = 3.25
$ws.Range("AV6").Value = 81
$ws.Range("AX6").Value = 7.5
$ws.Range("BA6").Value = 151
$ws.Range("BB6").Value = 251
$ws.Range("G6").Value = 1.6
$ws.Range("H6").Value = 3.6
$ws.Range("I6").Value = 6
$ws.Range("J6").Value = 2.3
$ws.Range("L6").Value = 7
$ws.Range("Q6").Value = 2.6
$ws.Range("R6").Value = 1.48
$ws.Range("Z6").Value = 11

# Row 7
$ws.Range("M7").Value = 1.14
$ws.Range("N7").Value = 5.5
$ws.Range("Q7").Value = 3.1
$ws.Range("R7").Value = 1.36

# Row 10
$ws.Range("AA10").Value = 21
$ws.Range("AC10").Value = 8
$ws.Range("AG10").Value = 351
$ws.Range("AH10").Value = 8.5
$ws.Range("AJ10").Value = 11
$ws.Range("AK10").Value = 29
$ws.Range("AL10").Value = 26
$ws.Range("AO10").Value = 15
$ws.Range("AP10").Value = 26
$ws.Range("AQ10").Value = 51
$ws.Range("AT10").Value = 2.5
$ws.Range("AV10").Value = 67
$ws.Range("AX10").Value = 4.75
$ws.Range("AY10").Value = 17
$ws.Range("BC10").Value = 251
$ws.Range("G10").Value = 2.4
$ws.Range("H10").Value = 3.2
$ws.Range("I10").Value = 3
$ws.Range("J10").Value = 3.2
$ws.Range("K10").Value = 2
$ws.Range("M10").Value = 1.08
$ws.Range("N10").Value = 8
$ws.Range("O10").Value = 1.4
$ws.Range("P10").Value = 2.75
$ws.Range("Q10").Value = 2.25
$ws.Range("R10").Value = 1.62
$ws.Range("S10").Value = 1.5
$ws.Range("T10").Value = 2.5
$ws.Range("U10").Value = 1.91
$ws.Range("V10").Value = 1.8
$ws.Range("X10").Value = 11
$ws.Range("Y10").Value = 10
$ws.Range("Z10").Value = 23

# Row 42
$ws.Range("AD42").Value = 7
$ws.Range("AI42").Value = 26
$ws.Range("AQ42").Value = 34
$ws.Range("AX42").Value = 7
$ws.Range("G42").Value = 1.65
$ws.Range("H42").Value = 3.5
$ws.Range("I42").Value = 5.5
$ws.Range("J42").Value = 2.38
$ws.Range("U42").Value = 2.38
$ws.Range("V42").Value = 1.53
$ws.Range("Z42").Value = 12
